$wb = $excel.ActiveWorkbook

# Sheet index 2 = "Percepciones" (1-based) per xl/workbook.xml sheet order:
# 1 Generales, 2 Percepciones, 3 Deducciones, 4 Otros Pagos
$wsPercepciones = $wb.Worksheets.Item("Percepciones")

# Add the new header cells AC3:AH3 on the "Percepciones" sheet
$wsPercepciones.Range("AC3").Value = "TP"
$wsPercepciones.Range("AD3").Value = "TD"
$wsPercepciones.Range("AE3").Value = "OP"
$wsPercepciones.Range("AF3").Value = "NETO TIMB"
$wsPercepciones.Range("AG3").Value = "NETO NOM"
$wsPercepciones.Range("AH3").Value = "-"

# Apply the grey fill style consistent with the rest of row 3 on that sheet
$wsPercepciones.Range("AC3:AH3").Interior.Color = 14474460

# Update the frozen-pane top left cell and selection to reflect scrolling to the new columns
$wsPercepciones.Activate()
$wsPercepciones.Range("J4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 10
$wsPercepciones.Range("AC3:AH3").Select()

# Make "Percepciones" the active (selected) sheet/tab of the workbook
$wsPercepciones.Activate()
